# Excel COM-interop script replicating the "data source files updates" commit.
#
# Summary of the edit:
#  1. pitstop sheet: add two new columns (H: tyre_before, I: tyre_after) with
#     per-row tyre-compound values.
#  2. Two brand-new worksheets appended at the end of the workbook:
#       - "weather"  : five weather readings (sky condition, temperature,
#                      humidity [as a %], wind speed, wind bearing).
#       - "altitude" : a single "delta" altitude reading.
#  3. Selection/active-cell bookkeeping that naturally falls out of the
#     above edits (pitstop -> lapsundercut -> weather -> altitude), ending
#     with "altitude" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. pitstop: tyre_before / tyre_after columns
# ---------------------------------------------------------------------------
$pitstop = $wb.Worksheets.Item("pitstop")

$pitstop.Cells.Item(1, 8).Value = "tyre_before"
$pitstop.Cells.Item(1, 9).Value = "tyre_after"

$tyreRows = @(
  @(2,4,2),
  @(3,4,2),
  @(4,2,3),
  @(5,4,2),
  @(6,2,4),
  @(7,3,2),
  @(8,4,2),
  @(9,2,4),
  @(10,3,2),
  @(11,3,2),
  @(12,2,4),
  @(13,3,2),
  @(14,2,4),
  @(15,3,2),
  @(16,2,4),
  @(17,4,3),
  @(18,3,3),
  @(19,3,4),
  @(20,2,3),
  @(21,3,4),
  @(22,2,3),
  @(23,3,4),
  @(24,3,2),
  @(25,2,3),
  @(26,3,2),
  @(27,4,2),
  @(28,3,2),
  @(29,2,4),
  @(30,3,2),
  @(31,2,3),
  @(32,3,2),
  @(33,3,2),
  @(34,2,4)
)

foreach ($r in $tyreRows) {
    $pitstop.Cells.Item($r[0], 8).Value = $r[1]
    $pitstop.Cells.Item($r[0], 9).Value = $r[2]
}

# column widths for the two new columns (best-fit to their content)
$pitstop.Columns.Item(8).ColumnWidth = 11.5703125
$pitstop.Columns.Item(9).ColumnWidth = 9.85546875

# leave the cursor where the author left it after typing the last value
$pitstop.Range("H35").Select()

# ---------------------------------------------------------------------------
# 2. lapsundercut: just a cursor move (no data change) - this sheet is no
#    longer the active tab afterwards
# ---------------------------------------------------------------------------
$lapsundercut = $wb.Worksheets.Item("lapsundercut")
$lapsundercut.Range("D57").Select()

# ---------------------------------------------------------------------------
# 3. New "weather" worksheet, appended after lapsundercut
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$weather = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$weather.Name = "weather"

$weather.Cells.Item(1, 1).Value = "Skycondition"
$weather.Cells.Item(1, 2).Value = "Clear"

$weather.Cells.Item(2, 1).Value = "Temperature"
$weather.Cells.Item(2, 2).Value = "69.97°F"

$weather.Cells.Item(3, 1).Value = "Humidity"
$weather.Cells.Item(3, 2).Value = 0.26
$weather.Cells.Item(3, 2).NumberFormat = "0%"

$weather.Cells.Item(4, 1).Value = "Wind speed"
$weather.Cells.Item(4, 2).Value = "7.92 mph"

$weather.Cells.Item(5, 1).Value = "Wind bearing"
$weather.Cells.Item(5, 2).Value = "165°"

$weather.Columns.Item(1).ColumnWidth = 12.85546875
$weather.Columns.Item(2).ColumnWidth = 9

$weather.Range("C2").Select()

# ---------------------------------------------------------------------------
# 4. New "altitude" worksheet, appended after weather - ends up the active
#    tab, matching the author's final saved state
# ---------------------------------------------------------------------------
$altitude = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $weather)
$altitude.Name = "altitude"

$altitude.Cells.Item(1, 1).Value = "delta"
$altitude.Cells.Item(1, 2).Value = 30.9

$altitude.Range("B2").Select()
